# Applies the "Update gh-pages to output generated at 456a3b4" edit.
# Workbook has 4 sheets: 1=展览 (Exhibition), 2=演出 (Performance),
# 3=本地生活 (Local Life), 4=全部类型 (All Types, a combined/sorted view).

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item(1)   # 展览
$wsShow = $wb.Worksheets.Item(2)   # 演出
$wsAll  = $wb.Worksheets.Item(4)   # 全部类型

# ---------------------------------------------------------------------
# 1) Sheet "展览" (Exhibition) - "想去人数" (F column) refresh for many rows
# ---------------------------------------------------------------------
$expoUpdates = @{
    2  = 2238
    3  = 120
    4  = 73
    5  = 702
    8  = 44
    9  = 2606
    10 = 1628
    11 = 1636
    12 = 314
    14 = 668
    15 = 837
    16 = 112
    17 = 337
    18 = 1097
    22 = 5726
    24 = 1007
    25 = 118
    28 = 260
    29 = 231
    31 = 1063
    32 = 838
    34 = 70
    36 = 426
    37 = 1197
    39 = 119
    40 = 197
    42 = 127
}
foreach ($row in $expoUpdates.Keys) {
    $wsExpo.Cells.Item($row, 6).Value = $expoUpdates[$row]
}

# ---------------------------------------------------------------------
# 2) Sheet "演出" (Performance) - rows 2 & 3 updates
# ---------------------------------------------------------------------
# Row2: "杭州·《卡农》永恒经典名曲音乐会" now shows as "已停售" instead of a 90 price
$wsShow.Range("G2").Value = "已停售"
# Row3: "杭州·才八点派对 Day1" interest count bumped 794 -> 795
$wsShow.Range("F3").Value = 795

# ---------------------------------------------------------------------
# 3) Sheet "全部类型" (All Types) - rows 2-5 re-sorted/refreshed content
# ---------------------------------------------------------------------
# Row 2: 杭州·才八点派对 Day1
$wsAll.Range("B2").Value = "2024.02.24"
$wsAll.Range("C2").Value = "杭州·才八点派对 Day1 "
$wsAll.Range("D2").Value = "通货路918号粮仓艺术公园7号楼 SoFunLivehouse"
$wsAll.Range("E2").Value = "2024.02.24 19:30-02.24 22:00"
$wsAll.Range("F2").Value = 795
$wsAll.Range("G2").Value = "不可售"
$wsAll.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=81692"
$wsAll.Range("I2").Value = "//i2.hdslb.com/bfs/openplatform/202402/TXDC8NrV1706866356112.jpeg"

# Row 3: 杭州·次元幻想动漫游戏嘉年华
$wsAll.Range("B3").Value = "2024.02.24"
$wsAll.Range("C3").Value = "杭州·次元幻想动漫游戏嘉年华"
$wsAll.Range("D3").Value = "德胜东路2539号 梦马汽车小镇"
$wsAll.Range("E3").Value = "2024.02.24 10:00-02.25 17:00"
$wsAll.Range("F3").Value = 2238
$wsAll.Range("G3").Value = 75
$wsAll.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=80425"
$wsAll.Range("I3").Value = "//i0.hdslb.com/bfs/openplatform/202401/ZlI1Z4Xh1704266621625.jpeg"

# Row 4: 杭州·大船文化·终极无伴奏——宁峰演绎伊萨伊与帕格尼尼音乐会
$wsAll.Range("B4").Value = "2024.02.25"
$wsAll.Range("C4").Value = "杭州·大船文化·终极无伴奏——宁峰演绎伊萨伊与帕格尼尼音乐会"
$wsAll.Range("D4").Value = "杭州市江干区新业路39号 杭州大剧院"
$wsAll.Range("E4").Value = "2024.02.25 19:30-02.25 21:10"
$wsAll.Range("F4").Value = 4
$wsAll.Range("G4").Value = 180
$wsAll.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=81311"
$wsAll.Range("I4").Value = "//i2.hdslb.com/bfs/openplatform/202401/k53cbfSX1706006394089.jpeg"

# Row 5: 杭州·才八点派对 Day2
$wsAll.Range("B5").Value = "2024.02.25"
$wsAll.Range("C5").Value = "杭州·才八点派对 Day2"
$wsAll.Range("D5").Value = "通货路918号粮仓艺术公园7号楼 SoFunLivehouse"
$wsAll.Range("E5").Value = "2024.02.25 19:30-02.25 22:00"
$wsAll.Range("F5").Value = 435
$wsAll.Range("G5").Value = "不可售"
$wsAll.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=81696"
$wsAll.Range("I5").Value = "//i2.hdslb.com/bfs/openplatform/202402/XaTpvowc1706867543851.jpeg"

# ---------------------------------------------------------------------
# 4) Sheet "全部类型" (All Types) - remaining "想去人数" (F column) refresh
# ---------------------------------------------------------------------
$allUpdates = @{
    6  = 120
    7  = 73
    8  = 702
    14 = 44
    15 = 2606
    16 = 1628
    17 = 1636
    18 = 314
    20 = 668
    22 = 837
    23 = 112
    24 = 337
    25 = 1097
    28 = 5726
    30 = 1007
    31 = 118
    34 = 260
    35 = 231
    37 = 1063
    38 = 838
    39 = 70
    40 = 426
    41 = 1198
    43 = 119
    44 = 197
    46 = 127
}
foreach ($row in $allUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allUpdates[$row]
}
